# Update the cryptos list with freshly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18 and 19 swap coin identity (Litecoin/ShibaInu order flipped) in
# addition to getting refreshed price/volume figures, so handle them
# explicitly together with everything else.
$updates = @(
    @{ Row = 2;  D = "27.825.73";    E = "  +1.38%  " },
    @{ Row = 3;  D = "1.886.09";     E = "  +1.15%  " },
    @{ Row = 4;  D = "1.013";        E = "  +1.26%  " },
    @{ Row = 5;  D = "334.81";       E = "  +1.51%  " },
    @{ Row = 6;  D = "1.011";        E = "  +1.12%  " },
    @{ Row = 7;  D = "0.4714";       E = "  -0.42%  " },
    @{ Row = 8;  D = "0.3919";       E = "  -1.31%  " },
    @{ Row = 9;  D = "47.87";        E = "  +1.27%  " },
    @{ Row = 10; D = "0.08047";      E = "  +0.27%  " },
    @{ Row = 11; D = "1.024";        E = "  +0.31%  " },
    @{ Row = 12; D = "22.04";        E = "  +2.01%  " },
    @{ Row = 13; D = "1.893.43";     E = "  +2.10%  " },
    @{ Row = 14; D = "5.973";        E = "  +0.13%  " },
    @{ Row = 15; D = "7.132";        E = "  -0.78%  " },
    @{ Row = 16; D = "1.013";        E = "  +1.12%  " },
    @{ Row = 17; D = "0.06737";      E = "  +2.85%  " },
    @{ Row = 18; B = "Litecoin";  C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "87.34";        E = "  +0.87%  " },
    @{ Row = 19; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib";      D = "0.00001049"; E = "  +0.86%  " },
    @{ Row = 20; D = "17.28";        E = "  -0.41%  " },
    @{ Row = 21; D = "1.010";        E = "  +0.85%  " },
    @{ Row = 22; D = "27.853.79";    E = "  +1.48%  " },
    @{ Row = 23; D = "5.511";        E = "  -0.15%  " },
    @{ Row = 24; D = "10.97";        E = "  -0.03%  " },
    @{ Row = 25; D = "2.336";        E = "  +1.50%  " },
    @{ Row = 26; D = "2.112.88";     E = "  +1.82%  " },
    @{ Row = 27; D = "159.13";       E = "  +3.21%  " },
    @{ Row = 28;                     E = "  -1.18%  " },
    @{ Row = 29; D = "2.101";        E = "  +0.54%  " },
    @{ Row = 30; D = "5.531";        E = "  -0.36%  " },
    @{ Row = 31; D = "121.71";       E = "  -0.63%  " },
    @{ Row = 32; D = "0.9758";       E = "  +1.61%  " },
    @{ Row = 33; D = "0.09479";      E = "  -0.48%  " },
    @{ Row = 34; D = "1.442";        E = "  -1.28%  " },
    @{ Row = 35; D = "3.640";        E = "  +1.48%  " },
    @{ Row = 36; D = "5.357";        E = "  +0.91%  " },
    @{ Row = 37; D = "0.06149";      E = "  +1.17%  " },
    @{ Row = 38; D = "0.02268";      E = "  +1.34%  " },
    @{ Row = 39; D = "1.219";        E = "  +0.25%  " },
    @{ Row = 40; D = "0.5999";       E = "  +0.39%  " },
    @{ Row = 41; D = "8.031";        E = "  -0.40%  " },
    @{ Row = 42; D = "0.1896";       E = "  -0.57%  " },
    @{ Row = 43; D = "10.30";        E = "  -0.57%  " },
    @{ Row = 44; D = "1.264";        E = "  -0.20%  " },
    @{ Row = 45; D = "0.5698";       E = "  +0.58%  " },
    @{ Row = 46; D = "12.28";        E = "  +0.61%  " },
    @{ Row = 47; D = "3.406";        E = "  -0.45%  " },
    @{ Row = 48; D = "1.937";        E = "  -0.17%  " },
    @{ Row = 49; D = "0.06921";      E = "  +2.05%  " },
    @{ Row = 50; D = "113.61";       E = "  +3.24%  " },
    @{ Row = 51;                     E = "  +5.81%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($u.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Column D holds price strings (e.g. "27.825.73", "1.013"). Many of
        # these look like plain numbers to Excel's auto-detection, so force
        # a text format while writing, then restore the cell to its
        # original (unstyled) state so only the text content changes.
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
